$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 202
$ws1.Range("F6").Value = 347
$ws1.Range("F8").Value = 2189
$ws1.Range("F9").Value = 374
$ws1.Range("F10").Value = 5369
$ws1.Range("F12").Value = 357

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 202
$ws4.Range("F7").Value = 347
$ws4.Range("F11").Value = 2189
$ws4.Range("F12").Value = 374
$ws4.Range("F13").Value = 5369
$ws4.Range("F15").Value = 357
